# Update the "Informe-03-030065-A-TM-TP" metadata sheet so that the
# header row uses human readable (capitalised) labels, the measure /
# dimension identifiers are corrected, and the "tipo" / "URI" rows line
# up correctly with their respective columns.
#
# Related issues fixed by this edit:
#   #8  Mejorar la generacion de SKOS Concept Schemes
#   #16 Incluir descripciones para algunas medidas en los DSDs
#   #17 Referenciada codelist que luego no tiene valores
#   #19 Anadir propiedad en el DSD que identifique el ambito territorial aplicable
#   #20 Generacion erronea de medidas en 01-080101-010105TC

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: human readable column headers
$ws.Range("A1").Value = "Segunda residencia"
$ws.Range("B1").Value = "Número hogares"
$ws.Range("C1").Value = "Provincia código"
$ws.Range("D1").Value = "Aragón"
$ws.Range("E1").Value = "Municipio código"
$ws.Range("F1").Value = "Provincia nombre"
$ws.Range("G1").Value = "Municipio nombre"

# Row 2: measure / dimension identifiers
$ws.Range("A2").Value = "iaest-measure:segunda-residencia"
$ws.Range("B2").Value = "iaest-measure:numero-hogares"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "sdmx-dimension:refArea"

# Row 3: "medida" / "dim" classification
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "dim"

# Row 4: data type / URI codelist references
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Provincia"
$ws.Range("G4").Value = "URI-Municipio"
